# Add a new "Save" column (H) to the sheet, mirroring the style used
# for the other header cells (B1:G1), and fill in the per-row save
# flag values for rows 2-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - set the text, then copy the formatting used by the
# other header cells (bold font, border, centered alignment) from G1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Per-row "Save" values
$saveValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
